# This weekly price sheet is organized as repeating pairs of rows
# ("Primera" / "Segunda" quality) going back in time as the row number
# increases. A new reporting week's prices were inserted at the top of the
# data block (rows 136/137), which pushes every older week down by one
# pair of rows (2 rows) - including the oldest pair, which now also gets
# appended again at the very bottom (rows 224/225), growing the sheet from
# 223 to 225 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right above the first "Primera"/"Segunda" pair of
# the data block. Excel shifts rows 136..223 down to 138..225 (values,
# styles/number formats and the sheet dimension all update automatically).
$ws.Rows("136:137").Insert()

# The freshly inserted rows 136/137 start out empty. Populate them with the
# same recurring "Primera"/"Segunda" pattern as the pair directly above
# them (rows 134/135) by copying that pair down into the new rows.
$ws.Rows("134:135").Copy()
$ws.Rows("136:136").PasteSpecial()

# Finally, stamp the new pair with the new reporting date (2023-06-13 as an
# Excel serial date number), which is the only thing that distinguishes
# this newest week's rows from the template pair they were copied from.
$ws.Cells.Item(136, 4).Value = 45090
$ws.Cells.Item(137, 4).Value = 45090
